# Updated input file conductor_1_coupling.xlsx to introduce a variable
# conductive heat transfer coefficient between solid components.
#
#  * add a new sheet "thermal_contact_resistance" right after "contact_HTC"
#    (and before "HTC_multiplier"), matching the other coupling-matrix
#    sheets in layout (header row, component-name row/column, 5x5 zeroed
#    matrix body).
#  * broaden the scope of "interf_thickness": it is no longer limited to
#    fluid elements, so its header text is updated accordingly.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Insert the new "thermal_contact_resistance" sheet
# ---------------------------------------------------------------------
$afterSheet = $wb.Worksheets.Item("contact_HTC")
$newSheet = $wb.Worksheets.Add($null, $afterSheet)
$newSheet.Name = "thermal_contact_resistance"

# Header / description cell
$newSheet.Range("C1").Value = "[m2K/W] thermal contact resistance between soldi component, should be >= 0; used if flag HTC_choice is 1"

# Component-name header row (mirrors every other coupling-matrix sheet)
$newSheet.Range("B2").Formula = "=contact_perimeter_flag!B2"
$newSheet.Range("C2").Formula = "=contact_perimeter_flag!C2"
$newSheet.Range("D2").Formula = "=contact_perimeter_flag!D2"
$newSheet.Range("E2").Formula = "=contact_perimeter_flag!E2"
$newSheet.Range("F2").Formula = "=contact_perimeter_flag!F2"

# Component-name header column
$newSheet.Range("A3").Formula = "=B2"
$newSheet.Range("A4").Formula = "=C2"
$newSheet.Range("A5").Formula = "=D2"
$newSheet.Range("A6").Formula = "=E2"
$newSheet.Range("A7").Formula = "=F2"

# Matrix body, all zeroed by default (same as the other coupling sheets)
$newSheet.Range("B3:F7").Value = 0

# ---------------------------------------------------------------------
# 2) Broaden the "interf_thickness" sheet description: it is now valid
#    for all conductor components (not just fluid elements), and the
#    matrix may no longer be symmetric.
# ---------------------------------------------------------------------
$interf = $wb.Worksheets.Item("interf_thickness")
$interf.Range("C1").Value = "[m] Thickness of the interface between conductor components. For interfaces between components that are not both channels, assing in cell (comp1,comp2) the thickess of the interface of comp1  when in contact with comp2 and in cell  (comp2,comp1) the thickess of the interface of comp2  when in contact with comp1. Therefore the matrix may no longer be symmetric and also the lower triangular region should be filled."

Write-Output "Edit complete"
